# Split the 'position' column (Variables sheet) into two columns:
# 'pivot' (figures / heading / stub) and 'order' (numeric order within the pivot group).
# Close #124

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Variables")

# Insert a brand-new column before the existing "position" column (column A).
# Everything that used to live in columns A:G (position, variable, en_long_name,
# type, en_note, en_elimination, en_domain) shifts right to B:H.
$ws.Columns.Item(1).Insert()

# New column A: "pivot" - replaces the old "position" header/values with the
# higher-level pivot role of each variable.
$ws.Range("A1").Value = "pivot"
$ws.Range("A2").Value = "figures"
$ws.Range("A3").Value = "heading"
$ws.Range("A4").Value = "stub"

# New column B: "order" - numeric ordering, only meaningful for heading/stub rows.
$ws.Range("B1").Value = "order"
$ws.Range("B3").Value = 1
$ws.Range("B4").Value = 1

# The old "type" column (now column E) used to carry "FIGURES" for the figures
# row; that information is now captured by the new A2 = "figures" value instead.
$ws.Range("E2").ClearContents()

# Update the active selections: the "Data" sheet no longer is the active tab,
# and its selection moves; the "Variables" sheet becomes the active tab with a
# new selection.
$wsData = $wb.Worksheets.Item("Data")
$wsData.Range("E20").Select()

$ws.Activate()
$ws.Range("F3").Select()
